$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> cells that change from 0 to 1
$changes = @{
    3  = @("G3", "H3")
    4  = @("D4", "E4")
    5  = @("H5")
    6  = @("D6", "E6")
    7  = @("H7")
    8  = @("H8")
    9  = @("H9")
    10 = @("H10")
    11 = @("H11")
    12 = @("D12", "E12")
    13 = @("H13")
    14 = @("G14", "H14")
    15 = @("H15")
    16 = @("H16")
    17 = @("H17")
    18 = @("H18")
}

foreach ($row in $changes.Keys) {
    foreach ($cellRef in $changes[$row]) {
        $ws.Range($cellRef).Value = 1
    }
}
